$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-28 Thursday" "2024-03-29 Friday"

Replace-Text "541×3=1623" "149×3=447"
Replace-Text "687×2=1374" "641×3=1923"
Replace-Text "704×5=3520" "282×3=846"
Replace-Text "640×9=5760" "671×3=2013"
Replace-Text "220×2=440" "436×2=872"

Replace-Text "220×8=1760" "863×2=1726"
Replace-Text "707×4=2828" "168×2=336"
Replace-Text "583×9=5247" "686×7=4802"
Replace-Text "780×4=3120" "945×9=8505"
Replace-Text "339×4=1356" "693×8=5544"

Replace-Text "461×5=2305" "478×8=3824"
Replace-Text "410×5=2050" "991×2=1982"
Replace-Text "253×7=1771" "380×2=760"
Replace-Text "578×5=2890" "979×9=8811"
Replace-Text "168×4=672" "297×8=2376"

Replace-Text "688×7=4816" "155×7=1085"
Replace-Text "302×8=2416" "392×9=3528"
Replace-Text "840×8=6720" "494×9=4446"
Replace-Text "208×6=1248" "550×7=3850"
Replace-Text "687×5=3435" "485×8=3880"

Replace-Text "766×9=6894" "578×8=4624"
Replace-Text "675×2=1350" "837×2=1674"
Replace-Text "332×8=2656" "554×2=1108"
Replace-Text "269×6=1614" "190×6=1140"
Replace-Text "328×3=984" "324×7=2268"
